$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BTAPOptions")

# ---------------------------------------------------------------------------
# 1) Cell VALUES, written in the exact order needed so newly-minted shared
#    strings land on the same ids as the target workbook.
# ---------------------------------------------------------------------------

$ws.Range("A230").Value = ":boiler_fuel"
$ws.Range("B230").Value = "Specify the fuel type of the primary and secondary boilers."

$ws.Range("C233").Value = "NaturalGasElecBackup"
$ws.Range("C235").Value = "ElectricityGasBackup"

$ws.Range("D231").Value = "Use NECB Default settings and take boiler fuel type from the primary heating fuel option."
$ws.Range("D232").Value = "Will force the creation of natural gas primary and secondary boilers and create baseboard hot-water heaters.  This will happen even if the NCEB default building would not normally inculde a boiler or hot-water baseboards."
$ws.Range("D233").Value = "Will force the creation of a natural gas primary and electric secondary boiler and create baseboard hot-water heaters.  This will happen even if the NCEB default building would not normally inculde a boiler or hot-water baseboards."
$ws.Range("D234").Value = "Will force the creation of electric primary and secondary boilers and create baseboard hot-water heaters.  This will happen even if the NCEB default building would not normally inculde a boiler or hot-water baseboards."
$ws.Range("D235").Value = "Will force the creation of an electric primary and natural gas secondary boiler and create baseboard hot-water heaters.  This will happen even if the NCEB default building would not normally inculde a boiler or hot-water baseboards."
$ws.Range("D236").Value = "Will force the creation of FuelOirNo2 primary and secondary boilers and create baseboard hot-water heaters.  This will happen even if the NCEB default building would not normally inculde a boiler or hot-water baseboards."

$ws.Range("A237").Value = ":boiler_cap_ratio"

$ws.Range("C239").Value = "0-0"

$ws.Range("D238").Value = "Use the default NECB behaviour.  If the boiler_fuel option was used, and it inculded disimilar primary and secondary boiler fuel types, the primary boiler will be 75% of the required boiler capacity and the secondary boiler will be 25% of the required boiler capacity."

$ws.Range("B237").Value = "Specify the ratio of primary and secondary boiler capacities.  These are based on the total required capacity for the boilers based on the sizing run.  They can sum to be more or less than 100%."

$ws.Range("D239").Value = "Use the default NECB behaviour regarless of boiler fuels."

$ws.Range("C240").Value = "80-20"

$ws.Range("D240").Value = "The first number sets the percent of the total required boiler capacity applied to the primary boiler.  The second number does the same for the secondary boiler.  In this case the primary boiler would have a capacity that is 80% of the total required boiler capacity and the secondary bolier's capacity would be 20% of the total required capacity.  Thes numbers can sum to more or less than 100."

# Remaining cells that reuse already-existing shared strings / are booleans.
$ws.Range("C231").Value = "NECB_Default"
$ws.Range("C232").Value = "NaturalGas"
$ws.Range("C234").Value = "Electricity"
$ws.Range("C236").Value = "FuelOilNo2"
$ws.Range("C238").Value = "NECB_Default"

$ws.Range("E231").Value = $true
$ws.Range("E232").Value = $true
$ws.Range("E233").Value = $true
$ws.Range("E234").Value = $true
$ws.Range("E235").Value = $true
$ws.Range("E236").Value = $true
$ws.Range("E238").Value = $true
$ws.Range("E239").Value = $true
$ws.Range("E240").Value = $true

# ---------------------------------------------------------------------------
# 2) Formatting
# ---------------------------------------------------------------------------

# Row 230 (":boiler_fuel" header) already carries the shaded / wrapped
# header style on A:H - only the text changed above. Copy that header
# format onto the two brand-new header rows (237 and the trailing blank
# separator row 241).
$ws.Range("A230:H230").Copy() | Out-Null
$ws.Range("A237:H237").PasteSpecial(-4122) | Out-Null
$ws.Range("A241:H241").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Rows.Item(230).RowHeight = 75
$ws.Rows.Item(237).RowHeight = 75
$ws.Rows.Item(241).RowHeight = 75

# Option (C) / option-description (D) columns: top-aligned, D wraps text -
# matches the styling used throughout the rest of the sheet for this kind
# of row.
$optionRows = 231,232,233,234,235,236,238,239,240
foreach ($r in $optionRows) {
    $ws.Cells.Item($r, 3).VerticalAlignment = -4160
    $ws.Cells.Item($r, 4).VerticalAlignment = -4160
    $ws.Cells.Item($r, 4).WrapText = $true
}

# Costing Included (E) column: boolean cells styled with the built-in
# "Good" cell style (green fill), like the other TRUE/FALSE flags already
# used with this workbook's conditional formatting blocks.
$ws.Range("E231").Style = "Good"
$ws.Range("E231").VerticalAlignment = -4160
$ws.Range("E231").Copy() | Out-Null
$ws.Range("E232:E236").PasteSpecial(-4122) | Out-Null
$ws.Range("E238:E240").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(231, 5).Value = $true
$ws.Cells.Item(232, 5).Value = $true
$ws.Cells.Item(233, 5).Value = $true
$ws.Cells.Item(234, 5).Value = $true
$ws.Cells.Item(235, 5).Value = $true
$ws.Cells.Item(236, 5).Value = $true
$ws.Cells.Item(238, 5).Value = $true
$ws.Cells.Item(239, 5).Value = $true
$ws.Cells.Item(240, 5).Value = $true

# Row heights for the wrapped description rows (the headless runtime does
# not auto-fit, so these are set explicitly to match the text-wrap driven
# heights Excel itself would have computed).
$ws.Rows.Item(231).RowHeight = 28.8
$ws.Rows.Item(232).RowHeight = 72
$ws.Rows.Item(233).RowHeight = 72
$ws.Rows.Item(234).RowHeight = 72
$ws.Rows.Item(235).RowHeight = 72
$ws.Rows.Item(236).RowHeight = 72
$ws.Rows.Item(238).RowHeight = 86.4
$ws.Rows.Item(239).RowHeight = 28.8
$ws.Rows.Item(240).RowHeight = 129.6

Write-Output "boiler_fuel / boiler_cap_ratio option blocks added"
